# Update "Jogos do Dia Betfair Back Lay" workbook for 2026-01-13
# - Remove the two trailing "English National League" / duplicate rows
#   (old rows 10 and 11), since the last fixture (Mexican Liga MX,
#   Puebla vs Mazatlan FC) has moved up into row 9 with refreshed odds.
# - Refresh the Betfair back/lay odds that changed across the remaining rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove old rows 10 and 11 (English National League / Southend vs
#     Gateshead, and the stale Mexican Liga MX row) -----------------------
$ws.Range("A10:A11").EntireRow.Delete()

# --- Row 2: Saudi Professional League, Al-Akhdoud vs Al-Kholood Club ----
$ws.Range("H2").Value = 2.44
$ws.Range("J2").Value = 3.3
$ws.Range("K2").Value = 3.8
$ws.Range("P2").Value = 1.87

# --- Row 3: German Bundesliga, Stuttgart vs Eintracht Frankfurt ---------
$ws.Range("N3").Value = 5.5
$ws.Range("P3").Value = 2.5
$ws.Range("R3").Value = 1.61
$ws.Range("S3").Value = 2.56

# --- Row 4: Saudi Professional League, Dhamk vs Al-Ittihad --------------
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 14
$ws.Range("J4").Value = 5.3
$ws.Range("P4").Value = 2.3
$ws.Range("Q4").Value = 1.61
$ws.Range("S4").Value = 2.54
$ws.Range("U4").Value = 1.83
$ws.Range("X4").Value = 26
$ws.Range("Y4").Value = 11.5
$ws.Range("AF4").Value = 120
$ws.Range("AG4").Value = 48
$ws.Range("AH4").Value = 34
$ws.Range("AI4").Value = 44
$ws.Range("AJ4").Value = 470
$ws.Range("AK4").Value = 200
$ws.Range("AL4").Value = 160
$ws.Range("AM4").Value = 180
$ws.Range("AN4").Value = 260

# --- Row 5: Saudi Professional League, Al-Fateh (KSA) vs Al Riyadh SC ---
$ws.Range("G5").Value = 2.1
$ws.Range("H5").Value = 4.1
$ws.Range("I5").Value = 5.5
$ws.Range("K5").Value = 5.1

# --- Row 6: German Bundesliga, Dortmund vs Werder Bremen ----------------
$ws.Range("P6").Value = 2.92
$ws.Range("R6").Value = 1.77
$ws.Range("T6").Value = 1.76
$ws.Range("AA6").Value = 310
$ws.Range("AD6").Value = 36
$ws.Range("AE6").Value = 140
$ws.Range("AM6").Value = 100

# --- Row 7: German Bundesliga, Hamburger SV vs Leverkusen ---------------
$ws.Range("G7").Value = 3.75
$ws.Range("Q7").Value = 1.69
$ws.Range("AC7").Value = 9.199999999999999

# --- Row 8: German Bundesliga, Mainz vs FC Heidenheim -------------------
$ws.Range("P8").Value = 1.92
$ws.Range("T8").Value = 1.99
$ws.Range("X8").Value = 14.5
$ws.Range("Y8").Value = 18.5
$ws.Range("Z8").Value = 50
$ws.Range("AB8").Value = 8.199999999999999
$ws.Range("AI8").Value = 110
$ws.Range("AJ8").Value = 17.5
$ws.Range("AK8").Value = 20
$ws.Range("AN8").Value = 11.5

# --- Row 9: now Mexican Liga MX, Puebla vs Mazatlan FC, with refreshed odds
$ws.Range("A9").Value = "Mexican Liga MX"
$ws.Range("C9").Value = "20:00:00"
$ws.Range("D9").Value = "Puebla"
$ws.Range("E9").Value = "Mazatlan FC"
$ws.Range("F9").Value = 1.95
$ws.Range("G9").Value = 2.36
$ws.Range("H9").Value = 3.35
$ws.Range("I9").Value = 4.7
$ws.Range("J9").Value = 3.35
$ws.Range("K9").Value = 5
$ws.Range("P9").Value = 1.91
$ws.Range("Q9").Value = 1.74
